# PROS-11454 CCJP KPI development
#
# Applies the workbook edits described in the commit:
#  - Bump the external-workbook index used by the "Validation_List" defined
#    name from [1] to [2] (the external link that feeds it was renumbered).
#  - Populate the "Include Stacking" (column B) cells for the new KPI rows
#    (rows 4, 5 and 7) with "N/A", matching the value already present in
#    sibling cells (e.g. C4) and elsewhere in the sheet (row 6).
#  - Tighten row 6's height slightly (32.95 -> 32.8).
#  - Move the active selection on the frozen bottom-right pane to C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Defined name: external reference index bump [1] -> [2] ---
$wb.Names.Item("Validation_List").RefersTo = "=[2]Set_up!`$A`$90:`$A`$124"

# --- New KPI data: fill previously-empty "Include Stacking" cells ---
$ws.Range("B4").Value = "N/A"
$ws.Range("B5").Value = "N/A"
$ws.Range("B7").Value = "N/A"

# --- Row height tweak for row 6 ---
$ws.Rows.Item(6).RowHeight = 32.8

# --- Selection moved to C10 ---
$ws.Range("C10").Select()
